$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b528e9603e9fd8d783f1bd4e210f295d5ff81a7/e2e/4e778f6c-ab8d-45bb-b14a-4043fa700c4e.md"
$targetDisplay = "4e778f6c-ab8d-45bb-b14a-4043fa700c4e.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99158f04da2d0235441d7759c7b7096804b55795/e2e/4e778f6c-ab8d-45bb-b14a-4043fa700c4e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b528e9603e9fd8d783f1bd4e210f295d5ff81a7/e2e/4e778f6c-ab8d-45bb-b14a-4043fa700c4e.md."

# --- zh-cn sheet: row 7 now has a completed handback ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $targetDisplay)
$wsZh.Range("I7").Font.Underline = $true
$wsZh.Range("I7").Font.Color = 15570276
$wsZh.Range("J7").Value = "4e778f6c-ab8d-45bb-b14a-4043fa700c4e.b5e6e9ae8d2243656bc3ddeb4237380534501d4a.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-22 15:03:49"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet: row 7 now has a completed handback ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $targetDisplay)
$wsDe.Range("I7").Font.Underline = $true
$wsDe.Range("I7").Font.Color = 15570276
$wsDe.Range("J7").Value = "4e778f6c-ab8d-45bb-b14a-4043fa700c4e.b5e6e9ae8d2243656bc3ddeb4237380534501d4a.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-22 15:03:56"
$wsDe.Range("P7").Value = $errorDetail
